$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Add the two new user rows under the existing data (rows 1-7 already used)
$ws.Range("A8").Value = "CB24190"
$ws.Range("B8").Value = "TEST11111"

$ws.Range("A9").Value = "CB21130"
$ws.Range("B9").Value = "test1234"

# Match the final selection state seen in the saved file
$ws.Range("B9").Select()
